$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.660.13"
$ws.Range("E2").Value = "  +7.30%  "
$ws.Range("D3").Value = "'1.743.70"
$ws.Range("E3").Value = "  +4.88%  "
$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("D5").Value = "'334.20"
$ws.Range("E5").Value = "  +5.62%  "
$ws.Range("D6").Value = "'0.9971"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").Value = "'0.3753"
$ws.Range("E7").Value = "  +3.30%  "
$ws.Range("E8").Value = "  +5.54%  "
$ws.Range("D9").Value = "'0.3436"
$ws.Range("E9").Value = "  +4.68%  "
$ws.Range("D10").Value = "'1.209"
$ws.Range("E10").Value = "  +5.78%  "
$ws.Range("D11").Value = "'0.07532"
$ws.Range("E11").Value = "  +6.25%  "
$ws.Range("D12").Value = "'0.9967"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "'6.498"
$ws.Range("E13").Value = "  +6.98%  "
$ws.Range("D14").Value = "'20.58"
$ws.Range("E14").Value = "  +4.44%  "
$ws.Range("D15").Value = "'7.099"
$ws.Range("E15").Value = "  +6.84%  "
$ws.Range("D16").Value = "'1.731.42"
$ws.Range("E16").Value = "  +3.85%  "
$ws.Range("D17").Value = "'0.00001098"
$ws.Range("E17").Value = "  +4.35%  "
$ws.Range("D18").Value = "'0.06702"
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("D19").Value = "'84.24"
$ws.Range("E19").Value = "  +5.58%  "
$ws.Range("D20").Value = "'0.9968"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").Value = "'16.84"
$ws.Range("E21").Value = "  +6.42%  "
$ws.Range("D22").Value = "'6.207"
$ws.Range("E22").Value = "  +4.32%  "
$ws.Range("D23").Value = "'13.19"
$ws.Range("E23").Value = "  +4.40%  "
$ws.Range("D24").Value = "'26.621.33"
$ws.Range("E24").Value = "  +7.07%  "
$ws.Range("D25").Value = "'2.478"
$ws.Range("E25").Value = "  +1.63%  "
$ws.Range("D26").Value = "'2.531"
$ws.Range("E26").Value = "  +4.25%  "
$ws.Range("D27").Value = "'1.419"
$ws.Range("E27").Value = "  +14.94%  "
$ws.Range("D28").Value = "'153.60"
$ws.Range("E28").Value = "  +3.29%  "
$ws.Range("D29").Value = "'19.70"
$ws.Range("E29").Value = "  +5.39%  "
$ws.Range("D30").Value = "'1.924.75"
$ws.Range("E30").Value = "  +3.98%  "
$ws.Range("D31").Value = "'132.60"
$ws.Range("E31").Value = "  +5.16%  "
$ws.Range("D32").Value = "'4.137"
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").Value = "'6.246"
$ws.Range("E33").Value = "  +6.39%  "
$ws.Range("D34").Value = "'0.08669"
$ws.Range("E34").Value = "  +2.19%  "
$ws.Range("D35").Value = "'1.724"
$ws.Range("E35").Value = "  +3.82%  "
$ws.Range("D36").Value = "'13.22"
$ws.Range("E36").Value = "  +7.32%  "
$ws.Range("D37").Value = "'5.501"
$ws.Range("E37").Value = "  +5.28%  "
$ws.Range("D38").Value = "'0.02376"
$ws.Range("E38").Value = "  +4.55%  "
$ws.Range("D39").Value = "'0.06389"
$ws.Range("E39").Value = "  +4.78%  "
$ws.Range("D40").Value = "'0.2196"
$ws.Range("E40").Value = "  +5.54%  "
$ws.Range("D41").Value = "'8.714"
$ws.Range("E41").Value = "  +4.39%  "
$ws.Range("D42").Value = "'1.246"
$ws.Range("E42").Value = "  -3.15%  "
$ws.Range("D43").Value = "'0.6308"
$ws.Range("E43").Value = "  +5.64%  "
$ws.Range("D44").Value = "'14.59"
$ws.Range("E44").Value = "  +13.17%  "
$ws.Range("D45").Value = "'0.9972"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").Value = "'3.910"
$ws.Range("E46").Value = "  +2.31%  "
$ws.Range("D47").Value = "'0.6132"
$ws.Range("E47").Value = "  +8.44%  "
$ws.Range("D48").Value = "'129.38"
$ws.Range("E48").Value = "  +2.65%  "
$ws.Range("D49").Value = "'2.081"
$ws.Range("E49").Value = "  +5.64%  "
$ws.Range("D50").Value = "'0.07354"
$ws.Range("E50").Value = "  +4.63%  "
$ws.Range("D51").Value = "'78.09"
$ws.Range("E51").Value = "  +3.86%  "
